$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77 (weekly update: new observation added at
# the top of the "Jengibre" block), pushing the former rows 77-87 down to
# 78-88 while preserving their data/formatting.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Range("A77").Value = 8
$ws.Range("B77").Value = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44918
$ws.Range("E77").Value = 4
$ws.Range("F77").Value = 100114007
$ws.Range("G77").Value = "Jengibre"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 450
$ws.Range("K77").Value = 14000
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 14500
$ws.Range("N77").Value = "`$/caja 13 kilos"
$ws.Range("O77").Value = "Perú"
$ws.Range("P77").Value = 1115
$ws.Range("Q77").Value = 13
$ws.Range("R77").Value = "Hortaliza"
